# Finished renaming LSI -> LIQ
# Rename every "TLSI_*" key in column A (rows 2-39 of the LIQ_dict sheet)
# to the equivalent "TLIQ_*" key. The German ("de") and English ("en")
# translation columns (B, C) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 39 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -like "TLSI_*") {
        $cell.Value = $val -replace "TLSI_", "TLIQ_"
    }
}

# A handful of rows ended up with their (previously customised) row height
# reset back to the sheet's default of 15, and row 28 to 30, once the
# rename was finished.
$ws.Rows.Item(24).RowHeight = 15
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(29).RowHeight = 15
$ws.Rows.Item(31).RowHeight = 15
$ws.Rows.Item(39).RowHeight = 15
